$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.86915842191171
$ws.Range("D2").Value = 9.927805132342959
$ws.Range("E2").Value = 16.42323057387616
$ws.Range("F2").Value = 31.71156053560778
$ws.Range("G2").Value = 3.646490667291693
$ws.Range("J2").Value = 11.65263026665763
$ws.Range("L2").Value = 13.32202308200447
$ws.Range("O2").Value = 24.00310828893458

$ws.Range("B3").Value = 20.35127642421058
$ws.Range("D3").Value = 9.936810372956389
$ws.Range("E3").Value = 16.332991680591
$ws.Range("F3").Value = 31.94508745417136
$ws.Range("G3").Value = 3.649161647005876
$ws.Range("J3").Value = 11.61126050081273
$ws.Range("L3").Value = 12.85844269249691
$ws.Range("O3").Value = 24.12866716332716

$ws.Range("B4").Value = 20.02719297182967
$ws.Range("D4").Value = 9.943590200616068
$ws.Range("E4").Value = 16.27987198549344
$ws.Range("F4").Value = 32.09993390950008
$ws.Range("G4").Value = 3.650887753069605
$ws.Range("J4").Value = 11.58738772418788
$ws.Range("L4").Value = 12.56421066642103
$ws.Range("O4").Value = 24.21409645666618

$ws.Range("B5").Value = 19.89375442074952
$ws.Range("D5").Value = 9.946668593510459
$ws.Range("E5").Value = 16.25881885471758
$ws.Range("F5").Value = 32.16590255813965
$ws.Range("G5").Value = 3.651612882755415
$ws.Range("J5").Value = 11.57805103455095
$ws.Range("L5").Value = 12.44203621937919
$ws.Range("O5").Value = 24.25099332752523

$ws.Range("B6").Value = 19.87151911505729
$ws.Range("D6").Value = 9.947198848443719
$ws.Range("E6").Value = 16.25535937327418
$ws.Range("F6").Value = 32.1770293772646
$ws.Range("G6").Value = 3.651734604430141
$ws.Range("J6").Value = 11.57652454742365
$ws.Range("L6").Value = 12.42161635493731
$ws.Range("O6").Value = 24.25724551392323

$ws.Range("B7").Value = 20.02539870647048
$ws.Range("D7").Value = 9.943630437805792
$ws.Range("E7").Value = 16.2795856284632
$ws.Range("F7").Value = 32.10081199429578
$ws.Range("G7").Value = 3.65089744435166
$ws.Range("J7").Value = 11.58726021108838
$ws.Range("L7").Value = 12.56257198536575
$ws.Range("O7").Value = 24.21458564027725

$ws.Range("B8").Value = 20.69195373677234
$ws.Range("D8").Value = 9.930651246822944
$ws.Range("E8").Value = 16.391652704011
$ws.Range("F8").Value = 31.78969347685013
$ws.Range("G8").Value = 3.64739379250754
$ws.Range("J8").Value = 11.63805323313779
$ws.Range("L8").Value = 13.16425048537036
$ws.Range("O8").Value = 24.04466357304409

$ws.Range("B9").Value = 21.94399522605746
$ws.Range("D9").Value = 9.915073166525755
$ws.Range("E9").Value = 16.62872802938846
$ws.Range("F9").Value = 31.2711657460739
$ws.Range("G9").Value = 3.641203060151174
$ws.Range("J9").Value = 11.74945128019825
$ws.Range("L9").Value = 14.26263648580427
$ws.Range("O9").Value = 23.77812868559212

$ws.Range("B10").Value = 22.82208904387286
$ws.Range("D10").Value = 9.909580963460767
$ws.Range("E10").Value = 16.81230357745236
$ws.Range("F10").Value = 30.94692491898685
$ws.Range("G10").Value = 3.63706452309414
$ws.Range("J10").Value = 11.83801982249923
$ws.Range("L10").Value = 15.01366249309113
$ws.Range("O10").Value = 23.6236840847175

$ws.Range("B11").Value = 23.21095275685897
$ws.Range("D11").Value = 9.908360148171374
$ws.Range("E11").Value = 16.89760112521165
$ws.Range("F11").Value = 30.81193299229924
$ws.Range("G11").Value = 3.635269775248225
$ws.Range("J11").Value = 11.8796675891556
$ws.Range("L11").Value = 15.34209737215308
$ws.Range("O11").Value = 23.56255894160656

$ws.Range("B12").Value = 23.35656910911404
$ws.Range("D12").Value = 9.908080242109472
$ws.Range("E12").Value = 16.93013610153041
$ws.Range("F12").Value = 30.76263035105858
$ws.Range("G12").Value = 3.634602713434545
$ws.Range("J12").Value = 11.89562450013497
$ws.Range("L12").Value = 15.46449248790118
$ws.Range("O12").Value = 23.54073766232912

$ws.Range("B13").Value = 23.3252826036442
$ws.Range("D13").Value = 9.908132434664449
$ws.Range("E13").Value = 16.92311903660751
$ws.Range("F13").Value = 30.77316751607986
$ws.Range("G13").Value = 3.634745819156432
$ws.Range("J13").Value = 11.89217978187581
$ws.Range("L13").Value = 15.4382215053186
$ws.Range("O13").Value = 23.54537812267493

$ws.Range("B14").Value = 23.22296615469427
$ws.Range("D14").Value = 9.908333472529135
$ws.Range("E14").Value = 16.90027322420151
$ws.Range("F14").Value = 30.80784034108094
$ws.Range("G14").Value = 3.635214644137251
$ws.Range("J14").Value = 11.88097670455544
$ws.Range("L14").Value = 15.3522068968918
$ws.Range("O14").Value = 23.56073705823162

$ws.Range("B15").Value = 23.16007781075184
$ws.Range("D15").Value = 9.908480326418159
$ws.Range("E15").Value = 16.88630939120002
$ws.Range("F15").Value = 30.82931544798359
$ws.Range("G15").Value = 3.635503447958974
$ws.Range("J15").Value = 11.87413841118455
$ws.Range("L15").Value = 15.29926103037098
$ws.Range("O15").Value = 23.57031782456217

$ws.Range("B16").Value = 22.79645275566948
$ws.Range("D16").Value = 9.909686342855283
$ws.Range("E16").Value = 16.80676324711184
$ws.Range("F16").Value = 30.95600011098164
$ws.Range("G16").Value = 3.6371835768428
$ws.Range("J16").Value = 11.83532456034602
$ws.Range("L16").Value = 14.99192591006861
$ws.Range("O16").Value = 23.62786358018202

$ws.Range("B17").Value = 22.5705838052025
$ws.Range("D17").Value = 9.910752490799496
$ws.Range("E17").Value = 16.75840632572092
$ws.Range("F17").Value = 31.03693266129103
$ws.Range("G17").Value = 3.638236744086415
$ws.Range("J17").Value = 11.81185458433843
$ws.Range("L17").Value = 14.79994462558052
$ws.Range("O17").Value = 23.66551390680601

$ws.Range("B18").Value = 22.43967986169726
$ws.Range("D18").Value = 9.91148604835033
$ws.Range("E18").Value = 16.73076237217673
$ws.Range("F18").Value = 31.08465907182226
$ws.Range("G18").Value = 3.638850775054409
$ws.Range("J18").Value = 11.79848362480374
$ws.Range("L18").Value = 14.68828251028451
$ws.Range("O18").Value = 23.68802825928355

$ws.Range("B19").Value = 22.39519167679131
$ws.Range("D19").Value = 9.911755126205065
$ws.Range("E19").Value = 16.72143242626539
$ws.Range("F19").Value = 31.10101989264142
$ws.Range("G19").Value = 3.639060099078623
$ws.Range("J19").Value = 11.79397877492096
$ws.Range("L19").Value = 14.65026523514092
$ws.Range("O19").Value = 23.69579836703318

$ws.Range("B20").Value = 22.59473128199386
$ws.Range("D20").Value = 9.910626552543516
$ws.Range("E20").Value = 16.76353660457169
$ws.Range("F20").Value = 31.02819540080565
$ws.Range("G20").Value = 3.638123776492304
$ws.Range("J20").Value = 11.81433978704905
$ws.Range("L20").Value = 14.82051023588282
$ws.Range("O20").Value = 23.66141698448043

$ws.Range("B21").Value = 23.25306430023911
$ws.Range("D21").Value = 9.908269483627443
$ws.Range("E21").Value = 16.90697740513833
$ws.Range("F21").Value = 30.79760666984506
$ws.Range("G21").Value = 3.635076598269341
$ws.Range("J21").Value = 11.88426235312073
$ws.Range("L21").Value = 15.37752566093082
$ws.Range("O21").Value = 23.55618969739877

$ws.Range("B22").Value = 23.67372942442434
$ws.Range("D22").Value = 9.907791465066062
$ws.Range("E22").Value = 17.00208139046217
$ws.Range("F22").Value = 30.65749612603216
$ws.Range("G22").Value = 3.633158328537294
$ws.Range("J22").Value = 11.93103931442992
$ws.Range("L22").Value = 15.7300209562626
$ws.Range("O22").Value = 23.49514895273578

$ws.Range("B23").Value = 23.4501259849907
$ws.Range("D23").Value = 9.907949827617459
$ws.Range("E23").Value = 16.95120590685054
$ws.Range("F23").Value = 30.73130073638579
$ws.Range("G23").Value = 3.634175466574864
$ws.Range("J23").Value = 11.9059780406051
$ws.Range("L23").Value = 15.54296678392881
$ws.Range("O23").Value = 23.52701617319286

$ws.Range("B24").Value = 22.58381746579538
$ws.Range("D24").Value = 9.910683113416379
$ws.Range("E24").Value = 16.76121671446108
$ws.Range("F24").Value = 31.03214178819595
$ws.Range("G24").Value = 3.638174822499168
$ws.Range("J24").Value = 11.81321584539383
$ws.Range("L24").Value = 14.81121653421597
$ws.Range("O24").Value = 23.66326649800766

$ws.Range("B25").Value = 21.61206772960554
$ws.Range("D25").Value = 9.918237064953768
$ws.Range("E25").Value = 16.56285698339748
$ws.Range("F25").Value = 31.40154606130078
$ws.Range("G25").Value = 3.64280551490047
$ws.Range("J25").Value = 11.71809950725734
$ws.Range("L25").Value = 13.97491692543571
$ws.Range("O25").Value = 23.84301862151885

